$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2  = @{ C = 1.392321641630434;  E = 1.710071460977503 }
    3  = @{ C = 1.004409005705997;  E = 1.642433761320072 }
    4  = @{ C = 2.418114148635109;  E = 2.828066716168021 }
    5  = @{ C = 1.89159218653383;   E = 2.544631191216329 }
    6  = @{ C = 1.554977796875501;  E = 1.312870290004287 }
    7  = @{ C = 0.6180254938795482; E = 0.7749619016293785 }
    8  = @{ C = 0.481899667566732;  E = 0.7487574275252262 }
    9  = @{ C = 1.905862317202089;  E = 1.389591155234515 }
    10 = @{ C = 2.671046044496239;  E = 2.125743999456575 }
    11 = @{ C = 1.917627847674064;  E = 2.694711744616662 }
    12 = @{ C = 1.119562422009102;  E = 1.831617848540201 }
    13 = @{ C = 1.344920716048192;  E = 1.037735724446631 }
    14 = @{ C = 2.195375580740766;  E = 1.872521508785896 }
    15 = @{ C = 2.542856270410665;  E = 2.961494745505977 }
    16 = @{ C = 0.3979826440748235; E = 2.008592810942544 }
    17 = @{ C = -2.604000402888396; E = -0.08252516517808228 }
    18 = @{ C = 1.122551915563408; E = 0.254631175783615 }
    19 = @{ C = 2.273132718878146; E = 1.620205313802381 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("E$row").Value = $vals.E
}
